$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "45.495.78"
$ws.Range("E2").Value = "  +3.02%  "
$ws.Range("D3").Value = "2.428.74"
$ws.Range("E3").Value = "  -0.53%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.72"
$ws.Range("E5").Value = "  +3.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.92"
$ws.Range("E6").Value = "  +4.57%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.531"
$ws.Range("E9").Value = "  +6.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.63"
$ws.Range("E10").Value = "  +0.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("E12").Value = "  -2.31%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.13"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.07"
$ws.Range("E14").Value = "  +1.66%  "
$ws.Range("D15").Value = "2.804.02"
$ws.Range("E15").Value = "  -0.03%  "
$ws.Range("D16").Value = "2.422.53"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "45.331.53"
$ws.Range("E18").Value = "  +2.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.25"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.35"
$ws.Range("E20").Value = "  -1.90%  "
$ws.Range("D21").Value = "0.0₃0923"
$ws.Range("E21").Value = "  +1.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.91"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "244.72"
$ws.Range("E23").Value = "  +1.94%  "
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("E25").Value = "  +0.54%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.64"
$ws.Range("E27").Value = "  +1.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -6.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.61"
$ws.Range("E29").Value = "  +1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "49.28"
$ws.Range("E30").Value = "  +2.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.03"
$ws.Range("E31").Value = "  +0.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.33"
$ws.Range("E32").Value = "  +9.33%  "
$ws.Range("E33").Value = "  +5.69%  "
$ws.Range("E34").Value = "  +0.92%  "
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  +0.86%  "
$ws.Range("E37").Value = "  -3.11%  "
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.88"
$ws.Range("E39").Value = "  -2.25%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "125.94"
$ws.Range("E40").Value = "  -4.67%  "
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("E42").Value = "  +0.55%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.49"
$ws.Range("E43").Value = "  -4.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0291"
$ws.Range("D45").Value = "1.926.09"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("E46").Value = "  -2.69%  "
$ws.Range("E47").Value = "  +1.73%  "
$ws.Range("E48").Value = "  +10.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.14"
$ws.Range("E49").Value = "  -2.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "76.72"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "53.89"
$ws.Range("E51").Value = "  +1.03%  "
